# "show tables and other interaction"
# Update the recorded timestamps in column I (HIS column) to plain date
# serials, and move the active selection from J21 to G11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column I: replace the fractional datetime serials with whole-day date
# serials (same numFmtId/style stays attached to the cell automatically).
$ws.Range("I2").Value = 42767
$ws.Range("I3").Value = 42768
$ws.Range("I4").Value = 42769
$ws.Range("I5").Value = 42770
$ws.Range("I6").Value = 42773
$ws.Range("I7").Value = 42774
$ws.Range("I8").Value = 42775
$ws.Range("I9").Value = 42783
$ws.Range("I10").Value = 42784
$ws.Range("I11").Value = 42785
$ws.Range("I12").Value = 42786
$ws.Range("I13").Value = 42787
$ws.Range("I14").Value = 42790
$ws.Range("I15").Value = 42791
$ws.Range("I16").Value = 42792
$ws.Range("I17").Value = 42795
$ws.Range("I18").Value = 42796
$ws.Range("I19").Value = 42797
$ws.Range("I20").Value = 42798

# Move the active cell / selection from J21 to G11.
$ws.Range("G11").Select()
